# Update the cryptos price/volume snapshot (D = Price, E = Volume(1h))
# on the active worksheet. Numeric-looking Price values are written with a
# leading apostrophe so Excel keeps them as literal text (preserving
# formatting such as trailing zeros) instead of auto-coercing to Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.848.53"
$ws.Range("D3").Value = "3.758.40"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'620.75"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").Value = "'181.77"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "3.757.00"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("D11").Value = "'6.32"
$ws.Range("E11").Value = "  -4.84%  "
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "'41.40"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").Value = "'0.0000261"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "4.382.09"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "3.762.93"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "69.956.98"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'7.62"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").Value = "'16.77"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "'508.83"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'9.63"
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'2.53"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").Value = "'87.42"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").Value = "'13.18"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D27").Value = "'11.15"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").Value = "'0.0000137"
$ws.Range("E28").Value = "  +23.04%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "'2.91"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").Value = "'7.92"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").Value = "'31.20"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  +4.75%  "
$ws.Range("D37").Value = "'6.22"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  +2.45%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("D41").Value = "'50.30"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("D42").Value = "'45.71"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("D43").Value = "'429.04"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'2.85"
$ws.Range("E45").Value = "  +1.97%  "
$ws.Range("D46").Value = "3.010.87"
$ws.Range("E46").Value = "  -3.88%  "
$ws.Range("D47").Value = "'0.0365"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "'27.52"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("D50").Value = "'137.02"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("D51").Value = "'2.51"
$ws.Range("E51").Value = "  +1.97%  "
